{"js": "// 1) Fix the instructions paragraph: merge all runs into one corrected sentence\n//    (\"-A partir del p,\" -> \"-A partir del 8,\") while keeping its sz/szCs (28) formatting.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph with the instructions text (\"-A partir del p, ...\")\nconst instructionsPara = paragraphs.items[2];\ninstructionsPara.load(\"text\");\nawait context.sync();\n\nif (instructionsPara.text.indexOf(\"-A partir del\") !== -1) {\n  const fixedText =\n    \"-A partir del 8, obten\u00e9 los n\u00fameros siguientes usando, en cada caso, \" +\n    \"una sola multiplicaci\u00f3n o divisi\u00f3n. Anot\u00e1 el c\u00e1lculo que realizaste. \" +\n    \"Luego comprob\u00e1 con calculadora:\";\n  instructionsPara.insertText(fixedText, \"Replace\");\n  await context.sync();\n}\n\n// 2) Turn the plain hyphen right before the final \"5\" into an en dash.\nconst hyphenHits = body.search(\"1500 - 5\", { matchCase: true });\nhyphenHits.load(\"items\");\nawait context.sync();\nif (hyphenHits.items.length > 0) {\n  hyphenHits.items[0].insertText(\"1500 \\u2013 5\", \"Replace\");\n  await context.sync();\n}\n\n// 3) Append two new blank paragraphs (Normal style, sz/szCs 48) after the\n//    numbers line, matching what Word produces for fresh Enter presses.\nparagraphs.load(\"items\");\nawait context.sync();\nconst numbersPara = paragraphs.items[paragraphs.items.length - 1];\n\nconst blankParagraphOoxml =\n  '<w:p><w:pPr><w:pStyle w:val=\"Normal\"/><w:rPr><w:sz w:val=\"48\"/><w:szCs w:val=\"48\"/></w:rPr></w:pPr></w:p>';\nconst packageXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  blankParagraphOoxml +\n  blankParagraphOoxml +\n  \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\nconst endRange = numbersPara.getRange(\"After\");\nendRange.insertOoxml(packageXml, \"After\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Fix the instructions paragraph: merge all its runs into one corrected\n#    sentence (\"-A partir del p,\" -> \"-A partir del 8,\") while keeping the\n#    paragraph's sz/szCs (28) run formatting.\n$instructionsPara = $d.Paragraphs(3)\n$instrRange = $instructionsPara.Range\n$instrRange.MoveEnd(1, -1) | Out-Null\n$instrRange.Text = \"-A partir del 8, obten\u00e9 los n\u00fameros siguientes usando, en cada caso, una sola multiplicaci\u00f3n o divisi\u00f3n. Anot\u00e1 el c\u00e1lculo que realizaste. Luego comprob\u00e1 con calculadora:\"\n\n# 2) Turn the plain hyphen right before the final \"5\" into an en dash.\n$hyphenRange = $d.Content\n$found = $hyphenRange.Find.Execute(\"1500 - 5\", $true)\nif ($found) {\n    $dash = [char]0x2013\n    $hyphenRange.Text = \"1500 $dash 5\"\n}\n\n# 3) Append two new blank paragraphs (Normal style, sz/szCs 48) after the\n#    numbers line, matching what Word produces for fresh Enter presses.\n$blankParagraphXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:pStyle w:val=\"Normal\"/><w:rPr><w:sz w:val=\"48\"/><w:szCs w:val=\"48\"/></w:rPr></w:pPr></w:p>'\n\n$numbersPara = $d.Paragraphs($d.Paragraphs.Count)\n$insertPoint1 = $d.Range($numbersPara.Range.End, $numbersPara.Range.End)\n$insertPoint1.InsertXML($blankParagraphXml)\n\n$firstBlankPara = $d.Paragraphs($d.Paragraphs.Count)\n$insertPoint2 = $d.Range($firstBlankPara.Range.End, $firstBlankPara.Range.End)\n$insertPoint2.InsertXML($blankParagraphXml)\n"}
